$wb = $excel.ActiveWorkbook

# "meta" sheet holds the key/value settings describing the chart.
$meta = $wb.Worksheets.Item("meta")

# Row 9 is currently just an empty placeholder row (A9 only carries the
# "key" style). Insert a fresh row above it so that placeholder moves down
# to become row 10, then fill the new row 9 with a "style" / "default"
# key-value pair (matching the other rows: column A uses the bold/orange
# "key" style, column B is plain).
$meta.Rows.Item(9).Insert()

$meta.Range("A9").Value = "style"
$meta.Range("B9").Value = "default"

# Copy the "key" cell formatting (style index 1) from an existing key cell
# onto the new A9 so it matches the rest of column A.
$meta.Range("A8").Copy()
$meta.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
